$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Add the new "Fig" paragraph style (linked to a new "Fig Char" character
#    style), based on Body Text, used for figure paragraphs.
# ---------------------------------------------------------------------------
$fig = $d.Styles.Add("Fig", 1)
$fig.BaseStyle = "BodyText"
$fig.LinkStyle = "FigChar"
$fig.ParagraphFormat.Alignment = 1
$fig.NoProofing = 1

$figChar = $d.Styles.Add("Fig Char", 2)
$figChar.BaseStyle = "BodyTextChar"
$figChar.LinkStyle = "Fig"
$figChar.Font.Name = "Ebrima"
$figChar.Font.Size = 11
$figChar.Font.SizeBi = 11
$figChar.NoProofing = 1

# ---------------------------------------------------------------------------
# 2. Re-style the (empty) paragraph that immediately follows the table: it
#    switches from "Body Text" to the new "Fig" style and is explicitly
#    left-aligned.
# ---------------------------------------------------------------------------
$table = $d.Tables.Item(1)
$afterTable = $table.Range.Next(4)
$targetPara = $afterTable.Paragraphs.Item(1)
$targetPara.Style = "Fig"
$targetPara.Alignment = 0
